# Agenda_points_for_jan_4th.docx — Jan 2 2022 update:
# Append two new bullet items to the last ("After Jan 6th but before Jan 9th:")
# list, right after "Fit the final version of the paper to the LaTeX format".
$d = $word.ActiveDocument

# The last paragraph in the document is the existing last bullet.
# Collapse its range to the end and insert a new paragraph there so the
# new paragraph inherits the same ListParagraph style / numbering (numId 2).
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.Collapse(0)            # wdCollapseEnd
$r.InsertParagraphAfter()

$newPara1 = $d.Paragraphs.Last.Range
$newPara1.Text = "Get complete pipeline on Github"

$r2 = $d.Paragraphs.Last.Range
$r2.Collapse(0)            # wdCollapseEnd
$r2.InsertParagraphAfter()

$newPara2 = $d.Paragraphs.Last.Range
$newPara2.Text = "Models?"
